# Update ticket-sold numbers for two events on the "展览" and "全部类型" sheets.
# 景德镇·第十四届瓷都ACG动漫游戏博览会: 1356 -> 1357
# 江西·广电·Unlimited Project 动漫游戏博览会: 3137 -> 3164
# 江西·2024南昌玛雅《次元之芯》主题动漫嘉年华: 32 -> 33

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F12").Value = 1357
$wsExpo.Range("F13").Value = 3164
$wsExpo.Range("F21").Value = 33

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F13").Value = 1357
$wsAll.Range("F14").Value = 3164
$wsAll.Range("F22").Value = 33
